$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data block: B118:B143 (with B118 header label, C140:C144 error calcs) ---

$ws.Range("B118").Value = "1-70-1 1e7"

$ws.Range("B119").Value = 0.085133
$ws.Range("B120").Value = 0.048994000000000003
$ws.Range("B121").Value = 0.141404
$ws.Range("B122").Value = 0.109572
$ws.Range("B123").Value = 0.057875999999999997
$ws.Range("B124").Value = 0.20765
$ws.Range("B125").Value = 0.064782000000000006
$ws.Range("B126").Value = 0.088372000000000006
$ws.Range("B127").Value = 0.019244000000000001
$ws.Range("B128").Value = 0.032551999999999998
$ws.Range("B129").Value = 0.044278999999999999
$ws.Range("B130").Value = 0.051378
$ws.Range("B131").Value = 0.056221
$ws.Range("B132").Value = 0.041706
$ws.Range("B133").Value = 0.051472999999999998
$ws.Range("B134").Value = 0.049985000000000002
$ws.Range("B135").Value = 0.048557999999999997
$ws.Range("B136").Value = 0.05296
$ws.Range("B137").Value = 0.026405000000000001
$ws.Range("B138").Value = 0.048238000000000003
$ws.Range("B139").Value = 0.067297999999999997
$ws.Range("B140").Value = -0.053821000000000001
$ws.Range("B141").Value = 0.073504
$ws.Range("B142").Value = -0.017589
$ws.Range("B143").Value = -0.13505200000000001

# Error-vs-actual formulas (mirrors the C column pattern used earlier in the sheet)
$ws.Range("C140").Formula = "=ABS(B140-`$A27)/`$A27"
$ws.Range("C141").Formula = "=ABS(B141-`$A28)/`$A28"
$ws.Range("C142").Formula = "=ABS(B142-`$A29)/`$A29"
$ws.Range("C143").Formula = "=ABS(B143-`$A30)/`$A30"
$ws.Range("C144").Formula = "=AVERAGE(C140:C143)"

$ws.Range("C140:C144").NumberFormat = "0%"

# --- View-state changes: scroll the window down to the new block and select C143 ---
$excel.ActiveWindow.ScrollRow = 112
$ws.Range("C143").Select()
